$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 348.92307
$ws.Range("I28").Value = 318.27274
$ws.Range("K28").Value = 318.27274
$ws.Range("M28").Value = 166.72726
$ws.Range("H40").Value = 1977.5714
$ws.Range("I40").Value = 1900
$ws.Range("J40").Value = 2008.6
$ws.Range("K40").Value = 1900
$ws.Range("L40").Value = 2008.6
$ws.Range("M40").Value = -1725
$ws.Range("N40").Value = -2358.6
$ws.Range("H43").Value = 25001500
$ws.Range("I43").Value = 25001500
$ws.Range("K43").Value = 25001500
$ws.Range("M43").Value = -25001431
$ws.Range("H92").Value = 282.33334
$ws.Range("I92").Value = 208.54546
$ws.Range("J92").Value = 398.2857
$ws.Range("K92").Value = 208.54546
$ws.Range("L92").Value = 398.2857
$ws.Range("M92").Value = 1039.45454
$ws.Range("N92").Value = -2894.2857
$ws.Range("H94").Value = 1490.6666
$ws.Range("I94").Value = 1490.6666
$ws.Range("K94").Value = 1490.6666
$ws.Range("M94").Value = -1039.6666
$ws.Range("H116").Value = 19214.715
$ws.Range("I116").Value = 37334.668
$ws.Range("J116").Value = 5624.75
$ws.Range("K116").Value = 37334.668
$ws.Range("L116").Value = 5624.75
$ws.Range("M116").Value = -33892.668
$ws.Range("N116").Value = -12508.75
$ws.Range("H137").Value = 2332540.2
$ws.Range("J137").Value = 5892206.5
$ws.Range("L137").Value = 17676619.5
$ws.Range("N137").Value = -17681719.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 879
$ws.Range("I5").Value = 879
$ws.Range("K5").Value = 879
$ws.Range("M5").Value = -767
$ws.Range("H32").Value = 349606.44
$ws.Range("I32").Value = 478098.9
$ws.Range("J32").Value = 12313.75
$ws.Range("K32").Value = 478098.9
$ws.Range("L32").Value = 12313.75
$ws.Range("M32").Value = -477811.9
$ws.Range("N32").Value = -12887.75
$ws.Range("H45").Value = 88078.086
$ws.Range("I45").Value = 168673.33
$ws.Range("K45").Value = 168673.33
$ws.Range("M45").Value = -168296.33
$ws.Range("H63").Value = 15487.424
$ws.Range("I63").Value = 4030.25
$ws.Range("K63").Value = 4030.25
$ws.Range("M63").Value = -3344.25
$ws.Range("H66").Value = 15487.424
$ws.Range("I66").Value = 4030.25
$ws.Range("K66").Value = 20151.25
$ws.Range("M66").Value = -16719.25
$ws.Range("H97").Value = 4245.3057
$ws.Range("I97").Value = 4399.5356
$ws.Range("K97").Value = 4399.5356
$ws.Range("M97").Value = -3903.5356
$ws.Range("H110").Value = 2280.9583
$ws.Range("I110").Value = 1429.1111
$ws.Range("J110").Value = 4836.5
$ws.Range("K110").Value = 1429.1111
$ws.Range("L110").Value = 4836.5
$ws.Range("M110").Value = 615.8888999999999
$ws.Range("N110").Value = -8926.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 879
$ws.Range("I4").Value = 879
$ws.Range("K4").Value = 879
$ws.Range("M4").Value = -764
$ws.Range("H20").Value = 796.3077
$ws.Range("I20").Value = 594.75
$ws.Range("J20").Value = 1118.8
$ws.Range("K20").Value = 594.75
$ws.Range("L20").Value = 1118.8
$ws.Range("M20").Value = -347.75
$ws.Range("N20").Value = -1612.8
$ws.Range("H22").Value = 1038.0667
$ws.Range("I22").Value = 728.53845
$ws.Range("K22").Value = 728.53845
$ws.Range("M22").Value = -555.53845
$ws.Range("H86").Value = 6158.074
$ws.Range("I86").Value = 4376.077
$ws.Range("J86").Value = 7812.7856
$ws.Range("K86").Value = 4376.077
$ws.Range("L86").Value = 7812.7856
$ws.Range("M86").Value = -3253.077
$ws.Range("N86").Value = -10058.7856
$ws.Range("H89").Value = 6158.074
$ws.Range("I89").Value = 4376.077
$ws.Range("J89").Value = 7812.7856
$ws.Range("K89").Value = 21880.385
$ws.Range("L89").Value = 39063.928
$ws.Range("M89").Value = -16264.385
$ws.Range("N89").Value = -50295.928

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2754.0652
$ws.Range("I31").Value = 2216.1316
$ws.Range("J31").Value = 5309.25
$ws.Range("K31").Value = 2216.1316
$ws.Range("L31").Value = 5309.25
$ws.Range("M31").Value = -1921.1316
$ws.Range("N31").Value = -5899.25
$ws.Range("H34").Value = 2754.0652
$ws.Range("I34").Value = 2216.1316
$ws.Range("J34").Value = 5309.25
$ws.Range("K34").Value = 2216.1316
$ws.Range("L34").Value = 5309.25
$ws.Range("M34").Value = -2014.1316
$ws.Range("N34").Value = -5713.25
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()
$ws.Range("H86").Value = 18083.084
$ws.Range("J86").Value = 28699.6
$ws.Range("L86").Value = 28699.6
$ws.Range("N86").Value = -30945.6
$ws.Range("H89").Value = 18083.084
$ws.Range("J89").Value = 28699.6
$ws.Range("L89").Value = 143498
$ws.Range("N89").Value = -154730

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 36171.43
$ws.Range("I41").Value = 1600
$ws.Range("K41").Value = 4800
$ws.Range("M41").Value = -4462
$ws.Range("H122").Value = 3922967
$ws.Range("J122").Value = 4222.5
$ws.Range("L122").Value = 38002.5
$ws.Range("N122").Value = -42902.5
$ws.Range("H130").Value = 6499.6665
$ws.Range("J130").Value = 5666.3335
$ws.Range("L130").Value = 16999.0005
$ws.Range("N130").Value = -27039.0005
$ws.Range("H138").Value = 3476.1
$ws.Range("I138").Value = 3471.8823
$ws.Range("K138").Value = 10415.6469
$ws.Range("M138").Value = -5275.6469
$ws.Range("H139").Value = 2726.4075
$ws.Range("I139").Value = 1845.5
$ws.Range("K139").Value = 5536.5
$ws.Range("M139").Value = -396.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 24976.715
$ws.Range("J95").Value = 24976.715
$ws.Range("L95").Value = 24976.715
$ws.Range("N95").Value = -30468.715
$ws.Range("H97").Value = 223097.56
$ws.Range("I97").Value = 200271.4
$ws.Range("J97").Value = 251630.25
$ws.Range("K97").Value = 200271.4
$ws.Range("L97").Value = 251630.25
$ws.Range("M97").Value = -199775.4
$ws.Range("N97").Value = -252622.25
$ws.Range("H102").Value = 50001852
$ws.Range("I102").Value = 55557390
$ws.Range("K102").Value = 55557390
$ws.Range("M102").Value = -55555768
$ws.Range("H104").Value = 65000
$ws.Range("J104").Value = 65000
$ws.Range("L104").Value = 65000
$ws.Range("N104").Value = -71988

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 10370.625
$ws.Range("I7").Value = 4994.2
$ws.Range("K7").Value = 4994.2
$ws.Range("M7").Value = -4882.2
$ws.Range("H22").Value = 3702.75
$ws.Range("J22").Value = 4205.9
$ws.Range("L22").Value = 4205.9
$ws.Range("N22").Value = -4795.9
$ws.Range("H27").Value = 3702.75
$ws.Range("J27").Value = 4205.9
$ws.Range("L27").Value = 4205.9
$ws.Range("N27").Value = -4419.9
$ws.Range("H46").Value = 7944.591
$ws.Range("I46").Value = 22759.6
$ws.Range("J46").Value = 3587.2354
$ws.Range("K46").Value = 22759.6
$ws.Range("L46").Value = 3587.2354
$ws.Range("M46").Value = -22571.6
$ws.Range("N46").Value = -3963.2354
$ws.Range("H61").Value = 23030.334
$ws.Range("I61").Value = 22045.5
$ws.Range("J61").Value = 25000
$ws.Range("K61").Value = 22045.5
$ws.Range("L61").Value = 25000
$ws.Range("M61").Value = -21843.5
$ws.Range("N61").Value = -25404
$ws.Range("H101").Value = 23563.166
$ws.Range("J101").Value = 23563.166
$ws.Range("L101").Value = 23563.166
$ws.Range("N101").Value = -30053.166
$ws.Range("H113").Value = 23030.334
$ws.Range("I113").Value = 22045.5
$ws.Range("J113").Value = 25000
$ws.Range("K113").Value = 22045.5
$ws.Range("L113").Value = 25000
$ws.Range("M113").Value = -19875.5
$ws.Range("N113").Value = -29340
$ws.Range("H126").Value = 10370.625
$ws.Range("I126").Value = 4994.2
$ws.Range("K126").Value = 14982.6
$ws.Range("M126").Value = -12512.6

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2859884
$ws.Range("I107").Value = 2490
$ws.Range("J107").Value = 3177372
$ws.Range("K107").Value = 7470
$ws.Range("L107").Value = 9532116
$ws.Range("M107").Value = -5550
$ws.Range("N107").Value = -9535956

Write-Host "Updated price/profit columns across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR."
